# Add the new "2022-Q4" sheet data for 600489-中金黄金, as described in the
# commit "feat: add 2022-Q4 data".
#
# 1. A brand-new worksheet "2022-Q4" is inserted right after "总计" (i.e.
#    right before "2022-Q3"), containing the per-fund holding table for the
#    new quarter.
# 2. The "总计" (totals) sheet gets a new row 2 for "2022-Q4" (count=23,
#    value=11.24), and every other existing row shifts down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet before "2022-Q3"
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

# Header row (B1:H1), bold + thin border + centered/top aligned, matching
# the look of the other quarterly sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q4.Cells.Item(1, 2 + $i)
    $cell.NumberFormat = "@"
    $cell.Value2 = $headers[$i]
}
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows (row 2 .. 24): code, name, size, position, weight, value, rank
$q4Data = @(
    @("004475", "华泰柏瑞富利灵活配置混合A", "34.73", "93.35", "9.13", "3.1708", 1),
    @("014597", "华泰柏瑞富利灵活配置混合C", "28.83", "93.35", "9.13", "2.6322", 1),
    @("003175", "华泰柏瑞多策略灵活配置混合A", "25.17", "92.05", "9.13", "2.2980", 1),
    @("015450", "华泰柏瑞多策略灵活配置混合C", "6.59", "92.05", "9.13", "0.6017", 1),
    @("002207", "前海开源金银珠宝主题精选混合C", "6.89", "90.36", "8.05", "0.5546", 5),
    @("161017", "富国中证500指数增强（LOF）", "68.12", "92.30", "0.76", "0.5177", 9),
    @("320006", "诺安灵活配置混合", "8.35", "77.61", "3.53", "0.2948", 9),
    @("001302", "前海开源金银珠宝主题精选混合A", "3.55", "90.36", "8.05", "0.2858", 5),
    @("003304", "前海开源沪港深核心资源灵活配置混合A", "3.30", "90.48", "7.93", "0.2617", 3),
    @("003305", "前海开源沪港深核心资源灵活配置混合C", "3.17", "90.48", "7.93", "0.2514", 3),
    @("003318", "景顺长城中证500行业中性低波动指数", "10.01", "93.81", "1.14", "0.1141", 2),
    @("011631", "西藏东财中证有色金属指数增强C", "1.69", "92.37", "3.88", "0.0656", 5),
    @("011630", "西藏东财中证有色金属指数增强A", "1.37", "92.37", "3.88", "0.0532", 5),
    @("160620", "鹏华中证A股资源产业指数（LOF）A", "1.72", "94.49", "2.23", "0.0384", 3),
    @("510170", "国联安上证大宗商品股票ETF", "1.72", "97.56", "2.10", "0.0361", 3),
    @("161715", "招商中证大宗商品股票指数（LOF）", "1.69", "94.56", "1.01", "0.0171", 9),
    @("013332", "富国中证500指数增强(LOF)C", "1.87", "92.30", "0.76", "0.0142", 9),
    @("512260", "华安中证500行业中性低波动ETF", "0.94", "97.66", "1.19", "0.0112", 2),
    @("000417", "国联安新精选灵活配置混合", "0.51", "57.88", "1.84", "0.0094", 8),
    @("012808", "鹏华中证A股资源产业指数（LOF）C", "0.24", "94.49", "2.23", "0.0054", 3),
    @("007943", "富安达中证 500 指数增强", "0.30", "78.45", "1.25", "0.0038", 4),
    @("002334", "汇丰晋信大盘波动精选股票A", "0.15", "90.32", "1.85", "0.0028", 3),
    @("002335", "汇丰晋信大盘波动精选股票C", "0.02", "90.32", "1.85", "0.0004", 3)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $row = 2 + $i
    $rec = $q4Data[$i]

    $idxCell = $q4.Cells.Item($row, 1)
    $idxCell.Value2 = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $codeCell = $q4.Cells.Item($row, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value2 = $rec[0]

    $nameCell = $q4.Cells.Item($row, 3)
    $nameCell.NumberFormat = "@"
    $nameCell.Value2 = $rec[1]

    $sizeCell = $q4.Cells.Item($row, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value2 = $rec[2]

    $posCell = $q4.Cells.Item($row, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value2 = $rec[3]

    $weightCell = $q4.Cells.Item($row, 6)
    $weightCell.NumberFormat = "@"
    $weightCell.Value2 = $rec[4]

    $valueCell = $q4.Cells.Item($row, 7)
    $valueCell.NumberFormat = "@"
    $valueCell.Value2 = $rec[5]

    $q4.Cells.Item($row, 8).Value2 = $rec[6]
}

$q4.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new "2022-Q4" row at the top of the
#    data (row 2), shifting the rest down by one row.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

for ($r = 9; $r -ge 2; $r--) {
    $destRow = $r + 1
    $totals.Cells.Item($destRow, 2).Value2 = $totals.Cells.Item($r, 2).Value2
    $totals.Cells.Item($destRow, 3).Value2 = $totals.Cells.Item($r, 3).Value2
    $totals.Cells.Item($destRow, 4).Value2 = $totals.Cells.Item($r, 4).Value2
}

$totals.Cells.Item(10, 1).Value2 = 8
$totals.Cells.Item(10, 1).Font.Bold = $true
$totals.Cells.Item(10, 1).HorizontalAlignment = -4108
$totals.Cells.Item(10, 1).VerticalAlignment = -4160
$totals.Cells.Item(10, 1).Borders.LineStyle = 1

$totals.Cells.Item(2, 2).Value2 = "2022-Q4"
$totals.Cells.Item(2, 3).Value2 = 23
$totals.Cells.Item(2, 4).Value2 = 11.24

$totals.Range("A1").Select()

# ---------------------------------------------------------------------
# Keep the original "2020-Q4" sheet (now last) as the active/selected tab,
# matching the source workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
